# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that preserves exact text formatting (no implicit
# numeric coercion / loss of trailing zeros, scientific notation, etc.).
# Cells in this sheet are plain-text (inlineStr) values such as "1.510" or
# "0.000008676" that must stay text, not become the Number 1.51 / 8.676E-06.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '29.220.33'
$ws.Range("D3").Value = '1.843.88'
$ws.Range("E3").Value = '  +0.13%  '
Set-TextValue "D4" '0.9995'
$ws.Range("E4").Value = '  +0.09%  '
Set-TextValue "D5" '242.77'
$ws.Range("E5").Value = '  +0.73%  '
Set-TextValue "D6" '0.6628'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.03%  '
Set-TextValue "D8" '44.88'
$ws.Range("E8").Value = '  +6.46%  '
Set-TextValue "D9" '0.07434'
$ws.Range("E9").Value = '  +0.14%  '
Set-TextValue "D10" '0.2954'
$ws.Range("E10").Value = '  -0.30%  '
Set-TextValue "D11" '23.28'
$ws.Range("E11").Value = '  +1.57%  '
Set-TextValue "D12" '0.07766'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '1.846.92'
$ws.Range("E13").Value = '  +1.70%  '
Set-TextValue "D14" '5.019'
$ws.Range("E14").Value = '  -0.41%  '
Set-TextValue "D15" '0.6719'
$ws.Range("E15").Value = '  -1.21%  '
Set-TextValue "D16" '83.44'
$ws.Range("E16").Value = '  -3.43%  '
Set-TextValue "D17" '6.187'
$ws.Range("E17").Value = '  -0.44%  '
Set-TextValue "D18" '0.000008676'
$ws.Range("E18").Value = '  +5.11%  '
$ws.Range("D19").Value = '29.238.78'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '2.101.08'
$ws.Range("E20").Value = '  +1.75%  '
Set-TextValue "D21" '12.54'
$ws.Range("E21").Value = '  -0.09%  '
Set-TextValue "D22" '226.73'
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("E23").Value = '  +0.10%  '
Set-TextValue "D24" '7.173'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("E25").Value = '  +0.05%  '
Set-TextValue "D26" '158.88'
$ws.Range("E26").Value = '  -0.88%  '
Set-TextValue "D27" '0.1408'
$ws.Range("E27").Value = '  -0.60%  '
Set-TextValue "D28" '8.631'
$ws.Range("E28").Value = '  -1.17%  '
Set-TextValue "D29" '18.04'
$ws.Range("E29").Value = '  -0.05%  '
Set-TextValue "D30" '1.510'
$ws.Range("E30").Value = '  +0.34%  '
Set-TextValue "D31" '4.136'
$ws.Range("E31").Value = '  -1.85%  '
Set-TextValue "D32" '4.053'
$ws.Range("E32").Value = '  -0.90%  '
Set-TextValue "D33" '1.190'
$ws.Range("E33").Value = '  -0.65%  '
Set-TextValue "D34" '0.05331'
$ws.Range("E34").Value = '  -0.43%  '
Set-TextValue "D35" '1.870'
$ws.Range("E35").Value = '  +0.07%  '
Set-TextValue "D36" '0.7468'
$ws.Range("E36").Value = '  -1.52%  '
Set-TextValue "D37" '1.156'
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '1.315.71'
$ws.Range("E39").Value = '  -1.31%  '
Set-TextValue "D40" '0.01802'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  +0.80%  '
Set-TextValue "D42" '6.389'
$ws.Range("E42").Value = '  +6.23%  '
Set-TextValue "D43" '0.9009'
$ws.Range("E43").Value = '  -2.11%  '
Set-TextValue "D44" '0.9996'
$ws.Range("E44").Value = '  -0.20%  '
Set-TextValue "D45" '103.52'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = '1.996.65'
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("B47").Value = 'XinFinNetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextValue "D47" '0.07860'
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D48" '65.35'
$ws.Range("E48").Value = '  +2.16%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D49" '0.00000000122'
$ws.Range("E49").Value = '  -1.55%  '
Set-TextValue "D50" '0.5144'
$ws.Range("E50").Value = '  -0.39%  '
Set-TextValue "D51" '1.753'
$ws.Range("E51").Value = '  -0.81%  '
